# Update TPM-derived metrics in the active sheet to reflect the new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 40.21121066666667
$ws.Range("N2").Value = 120.633632
$ws.Range("O2").Value = 0.3919993951244425
$ws.Range("P2").Value = 0.3919993951244425
$ws.Range("Q2").Value = 0.5681575992462222
$ws.Range("R2").Value = 5.113418393216
$ws.Range("S2").Value = 0.3919993951244425
$ws.Range("T2").Value = 0.3919993951244425

# Row 3
$ws.Range("O3").Value = 0.2238215523259795
$ws.Range("P3").Value = 0.2238215523259796
$ws.Range("S3").Value = 0.2238215523259795
$ws.Range("T3").Value = 0.2238215523259796

# Row 4
$ws.Range("M4").Value = 12.504156
$ws.Range("N4").Value = 37.512468
$ws.Range("O4").Value = 0.1218968916199506
$ws.Range("P4").Value = 0.1218968916199506
$ws.Range("Q4").Value = 0.176675388176
$ws.Range("R4").Value = 1.590078493584
$ws.Range("S4").Value = 0.1218968916199506
$ws.Range("T4").Value = 0.1218968916199506

# Row 5
$ws.Range("M5").Value = 26.90484566666667
$ws.Range("N5").Value = 80.71453700000001
$ws.Range("O5").Value = 0.2622821609296273
$ws.Range("P5").Value = 0.2622821609296273
$ws.Range("Q5").Value = 0.3801475327062223
$ws.Range("R5").Value = 3.421327794356
$ws.Range("S5").Value = 0.2622821609296273
$ws.Range("T5").Value = 0.2622821609296273
